$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.886.84"
$ws.Range("E2").Value = "  +0.85%  "
$ws.Range("D3").Value = "3.095.34"
$ws.Range("E3").Value = "  -0.37%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").Value = "540.21"
$ws.Range("D6").Value = "137.54"
$ws.Range("E6").Value = "  -0.63%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").Value = "3.088.56"
$ws.Range("E8").Value = "  -0.36%  "
$ws.Range("D9").Value = "0.498"
$ws.Range("E9").Value = "  +0.47%  "
$ws.Range("E10").Value = "  -2.51%  "
$ws.Range("D11").Value = "6.36"
$ws.Range("E11").Value = "  -3.46%  "
$ws.Range("D12").Value = "0.461"
$ws.Range("E12").Value = "  +0.69%  "
$ws.Range("D13").Value = "0.0000229"
$ws.Range("E13").Value = "  +4.64%  "
$ws.Range("D14").Value = "35.03"
$ws.Range("E14").Value = "  -0.49%  "
$ws.Range("D15").Value = "3.591.68"
$ws.Range("E15").Value = "  -0.54%  "
$ws.Range("D16").Value = "63.854.90"
$ws.Range("E16").Value = "  +0.73%  "
$ws.Range("E17").Value = "  +0.51%  "
$ws.Range("D18").Value = "3.093.30"
$ws.Range("E18").Value = "  -0.68%  "
$ws.Range("D19").Value = "6.73"
$ws.Range("E19").Value = "  +0.85%  "
$ws.Range("D20").Value = "490.38"
$ws.Range("E20").Value = "  -3.28%  "
$ws.Range("D21").Value = "13.55"
$ws.Range("E21").Value = "  -0.09%  "
$ws.Range("D22").Value = "0.705"
$ws.Range("E22").Value = "  -0.27%  "
$ws.Range("D23").Value = "7.23"
$ws.Range("E23").Value = "  -0.23%  "
$ws.Range("D24").Value = "80.14"
$ws.Range("E24").Value = "  +2.65%  "
$ws.Range("D25").Value = "12.29"
$ws.Range("E25").Value = "  -0.85%  "
$ws.Range("E26").Value = "  +0.10%  "
$ws.Range("D27").Value = "2.74"
$ws.Range("E27").Value = "  -0.89%  "
$ws.Range("D28").Value = "8.34"
$ws.Range("E28").Value = "  +0.65%  "
$ws.Range("D29").Value = "0.998"
$ws.Range("E29").Value = "  -0.15%  "
$ws.Range("D30").Value = "26.33"
$ws.Range("E30").Value = "  -0.19%  "
$ws.Range("E31").Value = "  -2.04%  "
$ws.Range("E32").Value = "  -0.40%  "
$ws.Range("D33").Value = "2.42"
$ws.Range("E33").Value = "  -5.06%  "
$ws.Range("D34").Value = "57.15"
$ws.Range("E34").Value = "  -1.45%  "
$ws.Range("D35").Value = "5.51"
$ws.Range("E35").Value = "  +5.28%  "
$ws.Range("D36").Value = "500.94"
$ws.Range("E36").Value = "  -5.49%  "
$ws.Range("D37").Value = "6.10"
$ws.Range("E37").Value = "  +1.86%  "
$ws.Range("D38").Value = "3.319.62"
$ws.Range("E38").Value = "  +7.72%  "
$ws.Range("D39").Value = "0.0401"
$ws.Range("E39").Value = "  -3.11%  "
$ws.Range("D40").Value = "0.0803"
$ws.Range("E40").Value = "  +1.12%  "
$ws.Range("D41").Value = "0.119"
$ws.Range("E41").Value = "  -2.08%  "
$ws.Range("D42").Value = "8.20"
$ws.Range("E42").Value = "  +0.85%  "
$ws.Range("D43").Value = "2.70"
$ws.Range("E43").Value = "  -1.57%  "
$ws.Range("D44").Value = "0.261"
$ws.Range("E44").Value = "  +2.97%  "
$ws.Range("D46").Value = "2.12"
$ws.Range("E46").Value = "  +3.22%  "
$ws.Range("D47").Value = "0.0₃0545"
$ws.Range("E47").Value = "  +7.05%  "
$ws.Range("D48").Value = "25.27"
$ws.Range("E48").Value = "  +3.76%  "
$ws.Range("D49").Value = "122.39"
$ws.Range("E49").Value = "  -0.43%  "
$ws.Range("E50").Value = "  +2.57%  "
$ws.Range("D51").Value = "2.37"
$ws.Range("E51").Value = "  -13.12%  "
